# "added parts list and orders"
#
# The parts list (A4:C35 — Qty / Part Number / description) gets sorted
# ascending by the "description" column (C), and the active selection
# moves from the old scroll position (B44 / topLeftCell A13) to B28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the parts table (A4:C35) by column C ascending, using the real
# Sort object so Excel records a <sortState> in the saved worksheet, the
# same as a manual Data > Sort would.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("C4:C35")) | Out-Null
$sortObj.SetRange($ws.Range("A4:C35"))
$sortObj.Header = 2
$sortObj.Orientation = 1
$sortObj.Apply()

# Update the visible selection to match the post-edit view.
$ws.Range("B28").Select() | Out-Null
